# Natmi following Dr Hou advice
# Re-computed LR-pair (Sema5a -> Plxnb3) edge-weight table: the expression-specificity
# statistics were recalculated (new Ligand/Receptor-expressing-cell counts etc.) and the
# "Target cluster" dimension now also includes the new "M2" cluster, so every
# Sending-cluster x Target-cluster combination (ECs/FAPs/sCs -> ECs/FAPs/sCs/M1/M2)
# is represented, growing the table from 9 data rows to 15 (A1:T10 -> A1:T16).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..J ligand stats, K..P receptor stats, Q..T edge stats.
$data = @(
    @('ECs', 'Sema5a', 'Plxnb3', 'ECs', 3, 1, 0.4898806666666666, 1.469642, 0.01965582386814743, 0.02319449744266509, 3, 1, 0.2951616666666667, 0.8854850000000001, 0.04824641865753034, 0.06350182003386917, 0.1445939940411111, 1.30134594637, 0.0009483231074013182, 0.001472892802380157),
    @('ECs', 'Sema5a', 'Plxnb3', 'FAPs', 3, 1, 0.4898806666666666, 1.469642, 0.01965582386814743, 0.02319449744266509, 2, 0.6666666666666666, 0.188123, 0.564369, 0.03075013472992963, 0.04047325326876763, 0.09215782065533333, 0.8294203858979999, 0.0006044192321733, 0.0009387567694387675),
    @('ECs', 'Sema5a', 'Plxnb3', 'M1', 3, 1, 0.4898806666666666, 1.469642, 0.01965582386814743, 0.02319449744266509, 3, 1, 0.7418906666666666, 2.225672, 0.1212676703798967, 0.1596122156766311, 0.3634378943804444, 3.270941049424, 0.002383615969887809, 0.003702125128329728),
    @('ECs', 'Sema5a', 'Plxnb3', 'M2', 3, 1, 0.4898806666666666, 1.469642, 0.01965582386814743, 0.02319449744266509, 2, 0.6666666666666666, 0.483482, 1.450446, 0.07902880902120335, 0.104017528090081, 0.2368484844813333, 2.131636360332, 0.001553376350630234, 0.002412634289277729),
    @('ECs', 'Sema5a', 'Plxnb3', 'sCs', 3, 1, 0.4898806666666666, 1.469642, 0.01965582386814743, 0.02319449744266509, 2, 1, 4.409136999999999, 8.818273999999999, 0.72070696721144, 0.632395182930651, 2.159950972984666, 12.959705837908, 0.01416608920805477, 0.01466808845323871),
    @('FAPs', 'Sema5a', 'Plxnb3', 'ECs', 3, 1, 13.025931, 39.077793, 0.5226485200912362, 0.6167418798615554, 3, 1, 0.2951616666666667, 0.8854850000000001, 0.04824641865753034, 0.06350182003386917, 3.844755503845001, 34.602799534605, 0.02521591931106044, 0.03916423186231865),
    @('FAPs', 'Sema5a', 'Plxnb3', 'FAPs', 3, 1, 13.025931, 39.077793, 0.5226485200912362, 0.6167418798615554, 2, 0.6666666666666666, 0.188123, 0.564369, 0.03075013472992963, 0.04047325326876763, 2.450477217513, 22.054294957617, 0.01607151240920385, 0.02496155030509259),
    @('FAPs', 'Sema5a', 'Plxnb3', 'M1', 3, 1, 13.025931, 39.077793, 0.5226485200912362, 0.6167418798615554, 3, 1, 0.7418906666666666, 2.225672, 0.1212676703798967, 0.1596122156766311, 9.663816633543998, 86.97434970189599, 0.06338036845896486, 0.09843953794527346),
    @('FAPs', 'Sema5a', 'Plxnb3', 'M2', 3, 1, 13.025931, 39.077793, 0.5226485200912362, 0.6167418798615554, 2, 0.6666666666666666, 0.483482, 1.450446, 0.07902880902120335, 0.104017528090081, 6.297803171741999, 56.680228545678, 0.04130429007950487, 0.06415196581282871),
    @('FAPs', 'Sema5a', 'Plxnb3', 'sCs', 3, 1, 13.025931, 39.077793, 0.5226485200912362, 0.6167418798615554, 2, 1, 4.409136999999999, 8.818273999999999, 0.72070696721144, 0.632395182930651, 57.43311433154699, 344.5986859892819, 0.3766764298325022, 0.3900245939360419),
    @('sCs', 'Sema5a', 'Plxnb3', 'ECs', 2, 1, 11.4071155, 22.814231, 0.4576956560406163, 0.3600636226957795, 3, 1, 0.2951616666666667, 0.8854850000000001, 0.04824641865753034, 0.06350182003386917, 3.366943222839167, 20.201659337035, 0.02208217623906858, 0.02286469536917036),
    @('sCs', 'Sema5a', 'Plxnb3', 'FAPs', 2, 1, 11.4071155, 22.814231, 0.4576956560406163, 0.3600636226957795, 2, 0.6666666666666666, 0.188123, 0.564369, 0.03075013472992963, 0.04047325326876763, 2.1459407892065, 12.875644735239, 0.01407420308855248, 0.01457294619423627),
    @('sCs', 'Sema5a', 'Plxnb3', 'M1', 2, 1, 11.4071155, 22.814231, 0.4576956560406163, 0.3600636226957795, 3, 1, 0.7418906666666666, 2.225672, 0.1212676703798967, 0.1596122156766311, 8.462832523038665, 50.776995138232, 0.05550368595104405, 0.05747055260302787),
    @('sCs', 'Sema5a', 'Plxnb3', 'M2', 2, 1, 11.4071155, 22.814231, 0.4576956560406163, 0.3600636226957795, 2, 0.6666666666666666, 0.483482, 1.450446, 0.07902880902120335, 0.104017528090081, 5.515135016171, 33.09081009702599, 0.03617114259106825, 0.03745292798797458),
    @('sCs', 'Sema5a', 'Plxnb3', 'sCs', 2, 1, 11.4071155, 22.814231, 0.4576956560406163, 0.3600636226957795, 2, 1, 4.409136999999999, 8.818273999999999, 0.72070696721144, 0.632395182930651, 50.29553501432349, 201.182140057294, 0.329864448170883, 0.2277025005413704)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2   # header is row 1, data starts at row 2
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Count; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $vals[$j]
    }
}
